$wb = $excel.ActiveWorkbook

# Sheet ALC, row 4
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 195
$ws.Range("I4").Value = 140
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 140
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = -26
$ws.Range("N4").Value = -478

# Sheet ALC, row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1031.75
$ws.Range("I32").Value = 759.4
$ws.Range("J32").Value = 1122.5333
$ws.Range("K32").Value = 759.4
$ws.Range("L32").Value = 1122.5333
$ws.Range("M32").Value = -433.4
$ws.Range("N32").Value = -1774.5333

# Sheet ALC, row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 313
$ws.Range("I38").Value = 298.5
$ws.Range("K38").Value = 895.5
$ws.Range("M38").Value = -523.5

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5722.222
$ws.Range("I40").Value = 2400
$ws.Range("K40").Value = 2400
$ws.Range("M40").Value = -2225

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2832.5
$ws.Range("I43").Value = 2915.3333
$ws.Range("J43").Value = 2749.6667
$ws.Range("K43").Value = 2915.3333
$ws.Range("L43").Value = 2749.6667
$ws.Range("M43").Value = -2846.3333
$ws.Range("N43").Value = -2887.6667

# Sheet ALC, row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1780.6875
$ws.Range("I58").Value = 500
$ws.Range("J58").Value = 2076.2307
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 6228.6921
$ws.Range("M58").Value = -1350
$ws.Range("N58").Value = -6528.6921

# Sheet ALC, row 61
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 63
$ws.Range("I61").Value = 63
$ws.Range("K61").Value = 189
$ws.Range("M61").Value = -17

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Sheet ALC, row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Sheet ALC, row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7439.7
$ws.Range("I113").Value = 4988.4443
$ws.Range("J113").Value = 9445.272000000001
$ws.Range("K113").Value = 4988.4443
$ws.Range("L113").Value = 9445.272000000001
$ws.Range("M113").Value = -1734.4443
$ws.Range("N113").Value = -15953.272

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3071.1428
$ws.Range("I122").Value = 1195.5
$ws.Range("K122").Value = 3586.5
$ws.Range("M122").Value = -1136.5

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2575.4546
$ws.Range("I132").Value = 2575.4546
$ws.Range("K132").Value = 7726.3638
$ws.Range("M132").Value = -5196.3638

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 832.4286
$ws.Range("I22").Value = 456.75
$ws.Range("J22").Value = 1333.3334
$ws.Range("K22").Value = 456.75
$ws.Range("L22").Value = 1333.3334
$ws.Range("M22").Value = -106.75
$ws.Range("N22").Value = -2033.3334

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 791.86664
$ws.Range("I122").Value = 816.5714
$ws.Range("J122").Value = 446
$ws.Range("K122").Value = 2449.7142
$ws.Range("L122").Value = 1338
$ws.Range("M122").Value = 0.2857999999996537
$ws.Range("N122").Value = -6238

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4316.6665
$ws.Range("I132").Value = 3750
$ws.Range("K132").Value = 11250
$ws.Range("M132").Value = -8720

# Sheet CUL, row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 499.4
$ws.Range("I2").Value = 236.75
$ws.Range("K2").Value = 1420.5
$ws.Range("M2").Value = -1307.5

# Sheet CUL, row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 313
$ws.Range("I7").Value = 210.25
$ws.Range("J7").Value = 518.5
$ws.Range("K7").Value = 630.75
$ws.Range("L7").Value = 1555.5
$ws.Range("M7").Value = -518.75
$ws.Range("N7").Value = -1779.5

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1719
$ws.Range("I68").Value = 1717
$ws.Range("J68").Value = 1719.8
$ws.Range("K68").Value = 5151
$ws.Range("L68").Value = 5159.4
$ws.Range("M68").Value = -4340
$ws.Range("N68").Value = -6781.4

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1719
$ws.Range("I71").Value = 1717
$ws.Range("J71").Value = 1719.8
$ws.Range("K71").Value = 15453
$ws.Range("L71").Value = 15478.2
$ws.Range("M71").Value = -11397
$ws.Range("N71").Value = -23590.2

# Sheet GSM, row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 148.26315
$ws.Range("I2").Value = 155.6923
$ws.Range("J2").Value = 132.16667
$ws.Range("K2").Value = 155.6923
$ws.Range("L2").Value = 132.16667
$ws.Range("M2").Value = -42.69229999999999
$ws.Range("N2").Value = -358.16667

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2819.2
$ws.Range("I80").Value = 2688
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2688
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1690
$ws.Range("N80").Value = -5996

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2819.2
$ws.Range("I83").Value = 2688
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 13440
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -8448
$ws.Range("N83").Value = -29984

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1193.6
$ws.Range("I122").Value = 1104.5
$ws.Range("J122").Value = 1550
$ws.Range("K122").Value = 3313.5
$ws.Range("L122").Value = 4650
$ws.Range("M122").Value = -863.5
$ws.Range("N122").Value = -9550

# Sheet LTW, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1849.6
$ws.Range("I22").Value = 1624.5
$ws.Range("J22").Value = 1999.6666
$ws.Range("K22").Value = 1624.5
$ws.Range("L22").Value = 1999.6666
$ws.Range("M22").Value = -1329.5
$ws.Range("N22").Value = -2589.6666

# Sheet LTW, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1849.6
$ws.Range("I27").Value = 1624.5
$ws.Range("J27").Value = 1999.6666
$ws.Range("K27").Value = 1624.5
$ws.Range("L27").Value = 1999.6666
$ws.Range("M27").Value = -1517.5
$ws.Range("N27").Value = -2213.6666

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3364.0356
$ws.Range("I46").Value = 2593.1875
$ws.Range("J46").Value = 4391.8335
$ws.Range("K46").Value = 2593.1875
$ws.Range("L46").Value = 4391.8335
$ws.Range("M46").Value = -2405.1875
$ws.Range("N46").Value = -4767.8335

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1100
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1100
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 2200
$ws.Range("N81").Value = -4322
$ws.Range("M81").ClearContents()

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1100
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1100
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 11000
$ws.Range("N84").Value = -21608
$ws.Range("M84").ClearContents()

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2817.1538
$ws.Range("I122").Value = 2802
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 8406
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -5956
$ws.Range("N122").Value = -13897
